$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (bold/bordered/centered style already applied to A1; extend to B1:G1)
$headers = @("group1", "group2", "meandiff", "p-adj", "lower", "upper", "reject")
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:G1").PasteSpecial(-4122) | Out-Null

# Data rows
$data = @(
    @(0, 3, 0.2922, 0.1927, -0.0903, 0.6746, $false),
    @(0, 5, 0.1666, 0.6418, -0.2159, 0.5491, $false),
    @(0, 6, -0.081, 0.9, -0.4635, 0.3014, $false),
    @(3, 5, -0.1256, 0.799, -0.508, 0.2569, $false),
    @(3, 6, -0.3732, 0.0584, -0.7557, 0.0093, $false),
    @(5, 6, -0.2476, 0.3271, -0.6301, 0.1348, $false)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
